# feat: add 2022-Q3 data
#
# Inserts a new worksheet "2022-Q3" (fund-holding snapshot) positioned
# right after "总计" and before the existing "2022-Q2" sheet, and updates
# the "总计" (totals) sheet with a new row for 2022-Q3 while shifting the
# existing 2022-Q2 / 2021-Q4 rows down.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)   # "总计"
$q2Sheet    = $wb.Worksheets.Item(2)   # "2022-Q2" (about to shift to position 3)

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet, inserted before the current Q2 sheet
#    so the final tab order is: 总计, 2022-Q3, 2022-Q2, 2021-Q4
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Add($q2Sheet)
$q3Sheet.Name = "2022-Q3"

# Headers (row 1) - bold, thin-bordered, centered/top, matching the other
# quarter sheets' header style.
$q3Sheet.Range("B1").Value = "基金代码"
$q3Sheet.Range("C1").Value = "基金名称"
$q3Sheet.Range("D1").Value = "基金规模"
$q3Sheet.Range("E1").Value = "股票总仓位"
$q3Sheet.Range("F1").Value = "仓位占比"
$q3Sheet.Range("G1").Value = "持有市值(亿元)"
$q3Sheet.Range("H1").Value = "仓位排名"

$headerRange = $q3Sheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1

# Data row (row 2) - keep the numeric-looking text columns as text so they
# render exactly like the source (e.g. no dropped trailing/leading zeros).
$q3Sheet.Range("B2:G2").NumberFormat = "@"

$q3Sheet.Range("A2").Value = 0
$q3Sheet.Range("B2").Value = "004685"
$q3Sheet.Range("C2").Value = "金元顺安元启灵活配置混合"
$q3Sheet.Range("D2").Value = "15.28"
$q3Sheet.Range("E2").Value = "77.14"
$q3Sheet.Range("F2").Value = "0.92"
$q3Sheet.Range("G2").Value = "0.1406"
$q3Sheet.Range("H2").Value = 2

# First-column marker cell (row 2) gets the same bold/border/center style.
$a2Range = $q3Sheet.Range("A2")
$a2Range.Font.Bold = $true
$a2Range.HorizontalAlignment = -4108  # xlCenter
$a2Range.VerticalAlignment = -4160    # xlTop
$a2Range.Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: existing rows shift down one, new row added
# ---------------------------------------------------------------------

# Row 2 now describes 2022-Q3
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.14

# Row 3 now describes 2022-Q2 (value unchanged from before)
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("C3").Value = 1
$totalSheet.Range("D3").Value = 0.1

# Row 4 (new) describes 2021-Q4 (value carried over from the old row 3)
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A4").PasteSpecial(-4122)  # xlPasteFormats

$totalSheet.Range("A4").Value = 2
$totalSheet.Range("B4").Value = "2021-Q4"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.05
